$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 8590.5
$ws.Range("I9").Value = 10196.8
$ws.Range("K9").Value = 10196.8
$ws.Range("M9").Value = -10027.8
$ws.Range("H27").Value = 1500
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 4500
$ws.Range("N27").Value = -4702
$ws.Range("H32").Value = 1200
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 900
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 900
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -1552
$ws.Range("H33").Value = 65865.12
$ws.Range("I33").Value = 111615.8
$ws.Range("J33").Value = 507
$ws.Range("K33").Value = 111615.8
$ws.Range("L33").Value = 507
$ws.Range("M33").Value = -111386.8
$ws.Range("N33").Value = -965
$ws.Range("H38").Value = 1115.5454
$ws.Range("I38").Value = 741.55554
$ws.Range("J38").Value = 2798.5
$ws.Range("K38").Value = 2224.66662
$ws.Range("L38").Value = 8395.5
$ws.Range("M38").Value = -1852.66662
$ws.Range("N38").Value = -9139.5
$ws.Range("H40").Value = 6540689
$ws.Range("I40").Value = 2510.3
$ws.Range("J40").Value = 15880945
$ws.Range("K40").Value = 2510.3
$ws.Range("L40").Value = 15880945
$ws.Range("M40").Value = -2335.3
$ws.Range("N40").Value = -15881295
$ws.Range("H43").Value = 7166.6665
$ws.Range("I43").Value = 7166.6665
$ws.Range("K43").Value = 7166.6665
$ws.Range("M43").Value = -7097.6665
$ws.Range("H58").Value = 3276.3333
$ws.Range("J58").Value = 8483
$ws.Range("L58").Value = 25449
$ws.Range("N58").Value = -25749
$ws.Range("H61").Value = 1428970
$ws.Range("J61").Value = 300
$ws.Range("L61").Value = 900
$ws.Range("N61").Value = -1244
$ws.Range("H74").Value = 1805518
$ws.Range("I74").Value = 1805518
$ws.Range("K74").Value = 1805518
$ws.Range("M74").Value = -1804582
$ws.Range("H76").Value = 3488
$ws.Range("J76").Value = 3275
$ws.Range("L76").Value = 3275
$ws.Range("N76").Value = -3905
$ws.Range("H77").Value = 1805518
$ws.Range("I77").Value = 1805518
$ws.Range("K77").Value = 9027590
$ws.Range("M77").Value = -9022910
$ws.Range("H79").Value = 3488
$ws.Range("J79").Value = 3275
$ws.Range("L79").Value = 3275
$ws.Range("N79").Value = -5459
$ws.Range("H106").Value = 2237.9443
$ws.Range("I106").Value = 2357
$ws.Range("K106").Value = 2357
$ws.Range("M106").Value = -1726
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H125").Value = 7888669.5
$ws.Range("J125").Value = 7814806
$ws.Range("L125").Value = 70333254
$ws.Range("N125").Value = -70338174
$ws.Range("H137").Value = 2214.7334
$ws.Range("I137").Value = 1993.5
$ws.Range("J137").Value = 3099.6667
$ws.Range("K137").Value = 5980.5
$ws.Range("L137").Value = 9299.000100000001
$ws.Range("M137").Value = -3430.5
$ws.Range("N137").Value = -14399.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 16500
$ws.Range("J31").Value = 32000
$ws.Range("L31").Value = 32000
$ws.Range("N31").Value = -32588
$ws.Range("H32").Value = 3451.8965
$ws.Range("I32").Value = 2503.75
$ws.Range("K32").Value = 2503.75
$ws.Range("M32").Value = -2216.75
$ws.Range("H101").Value = 37800
$ws.Range("J101").Value = 37800
$ws.Range("L101").Value = 37800
$ws.Range("N101").Value = -44290
$ws.Range("H122").Value = 5980.4116
$ws.Range("I122").Value = 5497.909
$ws.Range("K122").Value = 16493.727
$ws.Range("M122").Value = -14043.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3665.2942
$ws.Range("I86").Value = 3442.889
$ws.Range("K86").Value = 3442.889
$ws.Range("M86").Value = -2319.889
$ws.Range("H89").Value = 3665.2942
$ws.Range("I89").Value = 3442.889
$ws.Range("K89").Value = 17214.445
$ws.Range("M89").Value = -11598.445
$ws.Range("H107").Value = 211197.6
$ws.Range("I107").Value = 1994.5
$ws.Range("K107").Value = 1994.5
$ws.Range("M107").Value = -74.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 13293.818
$ws.Range("I86").Value = 9683.25
$ws.Range("J86").Value = 15357
$ws.Range("K86").Value = 9683.25
$ws.Range("L86").Value = 15357
$ws.Range("M86").Value = -8560.25
$ws.Range("N86").Value = -17603
$ws.Range("H89").Value = 13293.818
$ws.Range("I89").Value = 9683.25
$ws.Range("J89").Value = 15357
$ws.Range("K89").Value = 48416.25
$ws.Range("L89").Value = 76785
$ws.Range("M89").Value = -42800.25
$ws.Range("N89").Value = -88017
$ws.Range("H99").Value = 2913.4443
$ws.Range("I99").Value = 2902.625
$ws.Range("K99").Value = 2902.625
$ws.Range("M99").Value = -1404.625
$ws.Range("H126").Value = 2913.4443
$ws.Range("I126").Value = 2902.625
$ws.Range("K126").Value = 8707.875
$ws.Range("M126").Value = -6237.875
$ws.Range("H134").Value = 14708394
$ws.Range("I134").Value = 16668660
$ws.Range("K134").Value = 50005980
$ws.Range("M134").Value = -50003445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2843.4
$ws.Range("I17").Value = 2918.8572
$ws.Range("J17").Value = 2667.3333
$ws.Range("K17").Value = 8756.571599999999
$ws.Range("L17").Value = 8001.999899999999
$ws.Range("M17").Value = -8587.571599999999
$ws.Range("N17").Value = -8339.999899999999
$ws.Range("H129").Value = 3798.0527
$ws.Range("I129").Value = 1894.7142
$ws.Range("K129").Value = 5684.142599999999
$ws.Range("M129").Value = -684.1425999999992
$ws.Range("H131").Value = 1757.7273
$ws.Range("I131").Value = 1558.6666
$ws.Range("J131").Value = 1996.6
$ws.Range("K131").Value = 4675.9998
$ws.Range("L131").Value = 5989.799999999999
$ws.Range("M131").Value = 364.0002000000004
$ws.Range("N131").Value = -16069.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 440.8
$ws.Range("J3").Value = 52
$ws.Range("L3").Value = 52
$ws.Range("N3").Value = -284
$ws.Range("H43").Value = 1933.3334
$ws.Range("I43").Value = 1933.3334
$ws.Range("K43").Value = 1933.3334
$ws.Range("M43").Value = -1782.3334
$ws.Range("H46").Value = 39999
$ws.Range("J46").Value = 39999
$ws.Range("L46").Value = 39999
$ws.Range("N46").Value = -40311
$ws.Range("H97").Value = 2297.3076
$ws.Range("I97").Value = 1981.5
$ws.Range("J97").Value = 3350
$ws.Range("K97").Value = 1981.5
$ws.Range("L97").Value = 3350
$ws.Range("M97").Value = -1485.5
$ws.Range("N97").Value = -4342
$ws.Range("H107").Value = 6285
$ws.Range("I107").Value = 2663.3333
$ws.Range("K107").Value = 2663.3333
$ws.Range("M107").Value = -743.3332999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2481.5386
$ws.Range("I46").Value = 2423.6365
$ws.Range("K46").Value = 2423.6365
$ws.Range("M46").Value = -2235.6365
$ws.Range("H55").Value = 566.1
$ws.Range("J55").Value = 700.2
$ws.Range("L55").Value = 700.2
$ws.Range("N55").Value = -1046.2
$ws.Range("H68").Value = 2931.1667
$ws.Range("I68").Value = 2912.6667
$ws.Range("K68").Value = 2912.6667
$ws.Range("M68").Value = -2163.6667
$ws.Range("H71").Value = 2931.1667
$ws.Range("I71").Value = 2912.6667
$ws.Range("K71").Value = 14563.3335
$ws.Range("M71").Value = -10819.3335
$ws.Range("H122").Value = 14998.5
$ws.Range("I122").Value = 13997.667
$ws.Range("K122").Value = 41993.001
$ws.Range("M122").Value = -39543.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 46633.332
$ws.Range("J95").Value = 46633.332
$ws.Range("L95").Value = 46633.332
$ws.Range("N95").Value = -52125.332
$ws.Range("H113").Value = 870.86
$ws.Range("I113").Value = 880.1795
$ws.Range("K113").Value = 2640.5385
$ws.Range("M113").Value = -470.5384999999997
$ws.Range("H132").Value = 20837768
$ws.Range("I132").Value = 33335482
$ws.Range("J132").Value = 8243.666999999999
$ws.Range("K132").Value = 100006446
$ws.Range("L132").Value = 24731.001
$ws.Range("M132").Value = -100003916
$ws.Range("N132").Value = -29791.001

Write-Host "All updates applied"